$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Label" in column H, row 1 - reuse the existing header
# formatting (bold font, centered/top alignment, thin border) already used
# by B1:G1 so no new cell style gets minted.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Determine "Label" value per row: 0 for Control patients, 1 for MDD patients
# Patient order (rows 2-11, repeated for rows 12-21):
#   Control 26, Control 33, Control 36, Control 49, Control 2 -> Label 0
#   MDD 35, MDD 22, MDD 50, MDD 45, MDD 28 -> Label 1
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($block = 0; $block -lt 2; $block++) {
    for ($i = 0; $i -lt 10; $i++) {
        $row = 2 + ($block * 10) + $i
        $ws.Cells.Item($row, 8).Value = $labels[$i]
    }
}
